# Example6.xlsx "sample num 6 is cleaned up"
# The Aircraft_scheduling sheet had several E-column cells (departure/end
# time) that were computed with a formula like "=C4+0.9" (start time plus a
# duration, which can overflow past 24:00). Those formulas are replaced with
# plain literal time-of-day values (the wrapped/overflowed clock time), and
# the dependent F-column ("duration", a shared MOD(E-C,1) formula) is left
# as a formula so Excel recalculates its cached value from the new E value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aircraft_scheduling")

$ws.Range("E4").Value  = 0.26111111111111113
$ws.Range("E5").Value  = 0.2088888888888889
$ws.Range("E6").Value  = 0.011342592592592592
$ws.Range("E9").Value  = 0.17249999999999999
$ws.Range("E10").Value = 0.11763888888888889
$ws.Range("E11").Value = 0.45091435185185186
$ws.Range("E12").Value = 0.39869212962962958
$ws.Range("E13").Value = 0.56177083333333333
$ws.Range("E14").Value = 0.7220833333333333
$ws.Range("E15").Value = 0.58902777777777782
$ws.Range("E16").Value = 0.6759722222222222
$ws.Range("E17").Value = 0.64763888888888888
$ws.Range("E18").Value = 0.78541666666666676

# Move the active selection to E19, matching the saved view state.
$ws.Range("E19").Select()
